# Adds three new data rows (4, 5, 6) to the "Bulk Storage_Report_RNAseq_PBMC"
# sheet, matching the pre-existing HIV/HI/PBMC... row layout, and applies
# new gray/green highlight formatting to a couple of cells in the two
# newest rows.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bulk Storage_Report_RNAseq_PBMC")

# ---- Row 4: identical formatting to row 3, new Sample ID value ----
$ws.Range("A3:G3").Copy()
$ws.Range("A4:G4").PasteSpecial($xlPasteFormats)

$ws.Range("A4").Value = "HIV"
$ws.Range("B4").Value = "HI"
$ws.Range("C4").Value = "PBMC"
$ws.Range("D4").Value = "ATACseq"
$ws.Range("E4").Value = "STAS-06412_4"
$ws.Range("F4").Value = "Yellow"
$ws.Range("G4").Value = "AS13-08004"

# ---- Row 5: RNAseq entry, D/F highlighted (gray/green), G left-aligned ----
$ws.Range("A3:C3").Copy()
$ws.Range("A5:C5").PasteSpecial($xlPasteFormats)

$ws.Range("A3").Copy()
$ws.Range("E5").PasteSpecial($xlPasteFormats)

$ws.Range("A3").Copy()
$ws.Range("D5").PasteSpecial($xlPasteFormats)
$ws.Range("D5").Interior.Color = 12434877   # FFBDBDBD - gray highlight

$ws.Range("A3").Copy()
$ws.Range("F5").PasteSpecial($xlPasteFormats)
$ws.Range("F5").Interior.Color = 15332839   # FFE7F5E9 - green highlight

$ws.Range("G3").Copy()
$ws.Range("G5").PasteSpecial($xlPasteFormats)

$ws.Range("A5").Value = "HIV"
$ws.Range("B5").Value = "HI"
$ws.Range("C5").Value = "PBMC"
$ws.Range("D5").Value = "RNAseq"
$ws.Range("E5").Value = "FS07-06412_3"
$ws.Range("F5").Value = "Yellow"
$ws.Range("G5").Value = "FS07-06412"

# ---- Row 6: same formatting as row 5, different Sample IDs ----
$ws.Range("A3:C3").Copy()
$ws.Range("A6:C6").PasteSpecial($xlPasteFormats)

$ws.Range("A3").Copy()
$ws.Range("E6").PasteSpecial($xlPasteFormats)

$ws.Range("D5").Copy()
$ws.Range("D6").PasteSpecial($xlPasteFormats)

$ws.Range("F5").Copy()
$ws.Range("F6").PasteSpecial($xlPasteFormats)

$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial($xlPasteFormats)

$ws.Range("A6").Value = "HIV"
$ws.Range("B6").Value = "HI"
$ws.Range("C6").Value = "PBMC"
$ws.Range("D6").Value = "RNAseq"
$ws.Range("E6").Value = "FS07-05884_3"
$ws.Range("F6").Value = "Yellow"
$ws.Range("G6").Value = "FS07-05884"

$excel.CutCopyMode = $false
$ws.Range("A5:G6").Select()
